# JournalDeBord.xlsx - "Modfied logbook + documentation"
# Continue filling in the "Journal De Bord" log (rows 40-57) and fix a
# typo, plus add a hyperlink + page setup tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal De Bord")

# --- Fix existing rows -----------------------------------------------

# Row 40: fix typo "joural" -> "journal"
$ws.Range("C40").Value = "Question ajouter au journal de bord"

# Row 43: was "Implémentation des information personnel", now "Mot de passe "
$ws.Range("C43").Value = "Mot de passe "

# --- New journal rows 44-57 -------------------------------------------

# Row 44
$ws.Range("B44").NumberFormat = "h:mm"
$ws.Range("B44").Value = 0.39583333333333331
$ws.Range("C44").Value = "Implémentation des information personnel"

# Row 45
$ws.Range("B45").NumberFormat = "h:mm"
$ws.Range("B45").Value = 0.47916666666666669
$ws.Range("C45").Value = "Planning Effectif"

# Row 46 (new date: 4/6/2019)
$ws.Range("A46").NumberFormat = "m/d/yy"
$ws.Range("A46").Value = 43561
$ws.Range("B46").NumberFormat = "h:mm"
$ws.Range("B46").Value = 0.3125
$ws.Range("C46").Value = "Mail Experts"

# Row 47
$ws.Range("B47").NumberFormat = "h:mm"
$ws.Range("B47").Value = 0.31944444444444448
$ws.Range("C47").Value = "Implémentation des information personnel"

# Row 48
$ws.Range("B48").NumberFormat = "h:mm"
$ws.Range("B48").Value = 0.375
$ws.Range("C48").Value = "Test"

# Row 49 (no time value, only formatted)
$ws.Range("B49").NumberFormat = "h:mm"
$ws.Range("B49").Value = $null
$ws.Range("C49").Value = "Erreur "
$ws.Range("D49").Value = "requête SQL avec comme filtre le nouvelle email"

# Row 50 (no time value, only formatted)
$ws.Range("B50").NumberFormat = "h:mm"
$ws.Range("B50").Value = $null
$ws.Range("C50").Value = "Utiliser ~ pour REGEX en PHP"
$ws.Range("D50").Value = "https://www.sitepoint.com/community/t/check-whether-string-contains-numbers/5953"

# Row 51
$ws.Range("B51").NumberFormat = "h:mm"
$ws.Range("B51").Value = 0.52777777777777779
$ws.Range("C51").Value = "Implémentation Histoire"

# Row 52 - no B cell at all; D52 is a hyperlink
$ws.Range("C52").Value = "File input Bosstrap"
$ws.Hyperlinks.Add($ws.Range("D52"), "https://mdbootstrap.com/docs/jquery/forms/file-input/")

# Row 53 (no time value, only formatted)
$ws.Range("B53").NumberFormat = "h:mm"
$ws.Range("B53").Value = $null
$ws.Range("C53").Value = "Implémentation Fonction DB"

# Row 54
$ws.Range("B54").NumberFormat = "h:mm"
$ws.Range("B54").Value = 0.61805555555555558
$ws.Range("C54").Value = "RDV. avec Mme. Mota Stroppolo"

# Row 55 (no time value, only formatted)
$ws.Range("B55").NumberFormat = "h:mm"
$ws.Range("B55").Value = $null
$ws.Range("C55").Value = "Clément Christensen explication TRIM"

# Row 56
$ws.Range("B56").NumberFormat = "h:mm"
$ws.Range("B56").Value = 0.66666666666666663
$ws.Range("C56").Value = "Documentation"

# Row 57
$ws.Range("B57").NumberFormat = "h:mm"
$ws.Range("B57").Value = 0.68055555555555547
$ws.Range("C57").Value = "Planning Effectif"

# --- Column widths (best-effort; engine quantizes to whole char widths) --
$ws.Columns.Item(4).ColumnWidth = 81

# --- Page setup (A4 portrait) ------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Sheet2 ("Question") selection update ------------------------------
$ws2 = $wb.Worksheets.Item("Question")
$ws2.Range("B9").Select()

# --- Restore "Journal De Bord" as the active sheet/selection -----------
$ws.Activate()
$ws.Range("D49").Select()
